$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Copy formatting from the last fully-formatted data row (row 50) down
# into the new rows (51-58), so C/D get centered style, E gets the date
# style and G gets the centered style used by the "New" status column.
$ws.Range("A50:E50").Copy() | Out-Null
$ws.Range("A51:E58").PasteSpecial(-4122) | Out-Null
$ws.Range("G50").Copy() | Out-Null
$ws.Range("G51:G58").PasteSpecial(-4122) | Out-Null

# --- New requirement rows (User stories #31 - #35 plus related items)
$rows = @(
  @{ Row = 51; Id = 50; Text = "Pupil will get a notification when he got a feedback"; Priority = 4; Complexity = 2 },
  @{ Row = 52; Id = 51; Text = "Pupils can be connected to a class";                    Priority = 5; Complexity = 2 },
  @{ Row = 53; Id = 52; Text = "Teachers can be connected to a class";                  Priority = 5; Complexity = 2 },
  @{ Row = 54; Id = 53; Text = "Pupils can be disconnected from a class";               Priority = 4; Complexity = 2 },
  @{ Row = 55; Id = 54; Text = "Teachers can be disconnected from a class";             Priority = 4; Complexity = 2 },
  @{ Row = 56; Id = 55; Text = "A class can be created";                                Priority = 5; Complexity = 2 },
  @{ Row = 57; Id = 56; Text = "Class details can be edited";                           Priority = 4; Complexity = 2 },
  @{ Row = 58; Id = 57; Text = "A class can be removed";                                Priority = 4; Complexity = 2 }
)

foreach ($r in $rows) {
  $rowNum = $r.Row
  $ws.Cells.Item($rowNum, 1).Value = $r.Id
  $ws.Cells.Item($rowNum, 2).Value = $r.Text
  $ws.Cells.Item($rowNum, 3).Value = $r.Priority
  $ws.Cells.Item($rowNum, 4).Value = $r.Complexity
  $ws.Cells.Item($rowNum, 5).Value = 42686
  $ws.Cells.Item($rowNum, 7).Value = "New"
}

# --- Trailing blank id-only rows (keep the same centered "A" style)
$ws.Range("A50").Copy() | Out-Null
$ws.Range("A59:A61").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(59, 1).Value = 58
$ws.Cells.Item(60, 1).Value = 59
$ws.Cells.Item(61, 1).Value = 60

# --- Update the view to reflect the extended data: scroll down a bit and
# move the active selection to the new first empty row in column E.
$ws.Range("E62").Select() | Out-Null
